$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a phrase, then split off the trailing "(ACRONYM)" so that the
# surrounding "(" and ")" characters lose their bold formatting (matching the
# "Object/Function/Domain under Investigation" style cleanup).
# ---------------------------------------------------------------------------
function UnboldOuterParens($searchText) {
    $found = $d.Content
    $found.Find.Execute($searchText) | Out-Null
    if (-not $found.Find.Found) {
        return
    }
    $s = $found.Start
    $t = $found.Text
    $parenIdx = $t.LastIndexOf("(")
    $closeIdx = $t.LastIndexOf(")")

    $rParen = $d.Range($s + $parenIdx, $s + $parenIdx + 1)
    $rParen.Font.Bold = 0

    $rClose = $d.Range($s + $closeIdx, $s + $closeIdx + 1)
    $rClose.Font.Bold = 0
}

# 1. "Object under Investigation (OuI)" -- drop bold on the parentheses
UnboldOuterParens("Object under Investigation (OuI)")

# 2. "Function(s) under Investigation (FuI)" -- drop bold on the outer parens
UnboldOuterParens("Function(s) under Investigation (FuI)")

# 3. "Domain under Investigation (DuI):" -- drop the trailing colon, then
#    drop bold on the outer parens (same treatment as #1/#2 above).
$found = $d.Content
$found.Find.Execute("Domain under Investigation (DuI):") | Out-Null
if ($found.Find.Found) {
    $s = $found.Start
    $e = $found.End
    $colon = $d.Range($e - 1, $e)
    $colon.Delete()
}
UnboldOuterParens("Domain under Investigation (DuI)")

# 4. "Test criteria:" -> "Test criteria (TCR)"
#    Drop the colon, then append " (TCR)" where " " stays bold and "(TCR)"
#    is italic (not bold).
$found = $d.Content
$found.Find.Execute("Test criteria:") | Out-Null
if ($found.Find.Found) {
    $e = $found.End
    $colon = $d.Range($e - 1, $e)
    $colon.Delete()

    $insPoint = $d.Range($e - 1, $e - 1)
    $insPoint.InsertAfter(" (TCR)")

    $tcrStart = $e - 1
    $parenRange = $d.Range($tcrStart + 1, $tcrStart + 6)
    $parenRange.Font.Bold = 0
    $parenRange.Font.Italic = 1
}

# 5. Footer date update
$sec = $d.Sections.First
$ftr = $sec.Footers.Item(1)
$ftrRange = $ftr.Range
$ftrRange.Find.Execute("30-01-2020", $false, $false, $false, $false, $false, $true, 1, $false, "07-04-2021", 2) | Out-Null
